$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = 'Nf3'
$ws.Range("C13").Value = 'Nc6'
$ws.Range("H13").Value = 'Nf3'
$ws.Range("I13").Value = 'Nc6'
$ws.Range("C15").Value = 'Nd4'
$ws.Range("I15").Value = 'Bg4'
$ws.Range("B16").Value = 'Nxd4'
$ws.Range("H16").Value = 'Bg5'
$ws.Range("B17").Value = 'Qxd4'
$ws.Range("H17").Value = 'Be3'
$ws.Range("I17").Value = 'Qd7'
$ws.Range("B18").Value = 'Nc3'
$ws.Range("H18").Value = 'Nbd2'
$ws.Range("I19").Value = 'Bxf3'
$ws.Range("B20").Value = 'Qxe5+'
$ws.Range("C20").Value = 'Qe7'
$ws.Range("H20").Value = 'Qxf3'
$ws.Range("I20").Value = 'Nge7'
$ws.Range("B21").Value = 'Qxe7+'
$ws.Range("C21").Value = 'Nxe7'
$ws.Range("I21").Value = 'Qf7'
$ws.Range("B22").Value = 'Bd3'
$ws.Range("I22").Value = 'Nxc6'
$ws.Range("C23").Value = 'Nxd5'
$ws.Range("H23").Value = 'Bxa7'
$ws.Range("I23").Value = 'Nxa7'
$ws.Range("B24").Value = 'Re1+'
$ws.Range("C24").Value = 'Be7'
$ws.Range("B25").Value = 'Nxd5'
$ws.Range("I25").Value = 'Rxd5'
$ws.Range("C26").Value = 'Be6'
$ws.Range("H26").Value = 'Rd1'
$ws.Range("I26").Value = 'Rc5'
$ws.Range("B27").Value = 'Rxe6'
$ws.Range("H27").Value = 'Qd3'
$ws.Range("I27").Value = 'Bd6'
$ws.Range("B28").Value = 'Rxe7'
$ws.Range("C28").Value = 'Rhe8'
$ws.Range("H28").Value = 'Qf5+'
$ws.Range("I28").Value = 'Kb8'
$ws.Range("B29").Value = 'Rxg7'
$ws.Range("C29").Value = 'Re1+'
$ws.Range("I29").Value = 'Rxc4'
$ws.Range("B30").Value = 'Kh2'
$ws.Range("H30").Value = 'Bxc4'
$ws.Range("I30").Value = 'Qxc4'
$ws.Range("B31").Value = 'Bxf5+'
$ws.Range("C31").Value = 'Kb8'
$ws.Range("H31").Value = 'Ke2'
$ws.Range("I31").Value = 'Qxa3'
$ws.Range("B32").Value = 'Rg5'
$ws.Range("C32").Value = 'Rh8'
$ws.Range("H32").Value = 'Qe4'
$ws.Range("I32").Value = 'Qxb4'
$ws.Range("H33").Value = 'Qe3'
$ws.Range("I33").Value = 'Qb5+'
$ws.Range("B34").Value = 'Bb2'
$ws.Range("C34").Value = 'Rxa1'
$ws.Range("H34").Value = 'Kf3'
$ws.Range("I34").Value = 'Bc5'
$ws.Range("B35").Value = 'Bxa1'
$ws.Range("H35").Value = 'Qa3'
$ws.Range("I35").Value = 'Bxa3'
$ws.Range("H36").Value = 'Kg3'
$ws.Range("I36").Value = 'Qc3+'
$ws.Range("C37").Value = 'Rf8'
$ws.Range("I37").Value = 'Qxd2'
$ws.Range("C38").Value = 'Kc7'
$ws.Range("I38").Value = 'Qe3+'
$ws.Range("B39").Value = 'Be5+'
$ws.Range("C39").Value = 'Kb6'
$ws.Range("H39").Value = 'Kh4'
$ws.Range("I39").Value = 'Qxf4'
$ws.Range("B40").Value = 'Bd4+'
$ws.Range("C40").Value = 'Kc7'
$ws.Range("H40").Value = 'Kh5'
$ws.Range("B41").Value = 'Bg7'
$ws.Range("C41").Value = 'Rf7'
$ws.Range("H41").Value = 'Kh4'
$ws.Range("C42").Value = 'Kd6'
$ws.Range("H42").Value = 'Kh5'
$ws.Range("I42").Value = 'Qg3'
$ws.Range("B43").Value = 'Ba1'
$ws.Range("C43").Value = 'Kxd5'
$ws.Range("H43").Value = 'Rd1'
$ws.Range("I43").Value = 'Qxh3#'
$ws.Range("B44").Value = 'Rh5'
$ws.Range("C44").Value = 'Re7'
# Clear contents of cells whose chess-move data was removed in the diff
$ws.Range("B45:C47").ClearContents()
